$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.151.15"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").Value = "3.805.25"
$ws.Range("E3").Value = "  +7.13%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("B5").Value = "Solana"
$ws.Range("C5").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D5").Value = "'141.84"
$ws.Range("E5").Value = "  +7.60%  "
$ws.Range("B6").Value = "BNB"
$ws.Range("C6").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'416.13"
$ws.Range("E6").Value = "  -0.63%  "
$ws.Range("D7").Value = "3.790.51"
$ws.Range("E7").Value = "  +7.07%  "
$ws.Range("D8").Value = "'0.639"
$ws.Range("E8").Value = "  -3.16%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'0.759"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").Value = "'0.178"
$ws.Range("E11").Value = "  +5.07%  "
$ws.Range("D12").Value = "'0.0000372"
$ws.Range("E12").Value = "  +28.12%  "
$ws.Range("D13").Value = "'43.52"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "'10.29"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").Value = "4.378.76"
$ws.Range("E15").Value = "  +6.32%  "
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "3.816.36"
$ws.Range("E17").Value = "  +7.55%  "
$ws.Range("D18").Value = "'20.65"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "'13.33"
$ws.Range("E19").Value = "  +4.21%  "
$ws.Range("D20").Value = "'1.12"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").Value = "67.365.48"
$ws.Range("E21").Value = "  +1.91%  "
$ws.Range("D22").Value = "'437.31"
$ws.Range("E22").Value = "  -2.61%  "
$ws.Range("D23").Value = "'15.20"
$ws.Range("E23").Value = "  +15.19%  "
$ws.Range("D24").Value = "'89.22"
$ws.Range("E24").Value = "  -1.37%  "
$ws.Range("D25").Value = "'3.09"
$ws.Range("E25").Value = "  -5.14%  "
$ws.Range("D26").Value = "'37.45"
$ws.Range("E26").Value = "  +8.86%  "
$ws.Range("D27").Value = "'3.29"
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("D28").Value = "'9.75"
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").Value = "'5.17"
$ws.Range("E29").Value = "  +7.10%  "
$ws.Range("D30").Value = "'12.60"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "'0.122"
$ws.Range("E31").Value = "  +3.55%  "
$ws.Range("D32").Value = "'2.75"
$ws.Range("E32").Value = "  -1.89%  "
$ws.Range("D33").Value = "'7.21"
$ws.Range("E33").Value = "  -1.47%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").Value = "'41.85"
$ws.Range("E34").Value = "  +6.80%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "'0.161"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "'57.66"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "'0.0482"
$ws.Range("E38").Value = "  -5.02%  "
$ws.Range("E39").Value = "  +27.68%  "
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("D41").Value = "0.0₃0680"
$ws.Range("E41").Value = "  -7.88%  "
$ws.Range("D42").Value = "'0.996"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "'3.40"
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("D44").Value = "'147.81"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").Value = "'3.19"
$ws.Range("E45").Value = "  +23.81%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'4.41"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D47").Value = "'2.10"
$ws.Range("E47").Value = "  +4.64%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'26.24"
$ws.Range("E48").Value = "  +19.55%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.85"
$ws.Range("E49").Value = "  -6.34%  "
$ws.Range("D50").Value = "'2.57"
$ws.Range("E50").Value = "  -7.60%  "
$ws.Range("D51").Value = "'0.300"
$ws.Range("E51").Value = "  -3.77%  "
